# Apply targeted odds updates to Sheet1 for the 2026-01-06 Betfair Back/Lay workbook.
# The workbook already has the target rows/columns populated; we only need to
# overwrite the specific cell values that changed, matching the upstream diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 1.91
$ws.Range("J2").Value = 3.7
$ws.Range("K2").Value = 3.75
$ws.Range("L2").Value = 1.45
$ws.Range("N2").Value = 3.55
$ws.Range("P2").Value = 1.85
$ws.Range("S2").Value = 3.95
$ws.Range("T2").Value = 1.97
$ws.Range("X2").Value = 13.5
$ws.Range("AA2").Value = 120
$ws.Range("AC2").Value = 8.199999999999999
$ws.Range("AD2").Value = 19
$ws.Range("AG2").Value = 9.800000000000001
$ws.Range("AJ2").Value = 21
$ws.Range("F3").Value = 1.06
$ws.Range("M3").Value = 1.03
$ws.Range("N3").Value = 1.26
$ws.Range("O3").Value = 1.03
$ws.Range("Q3").Value = 1.06
$ws.Range("F4").Value = 1.77
$ws.Range("G4").Value = 1.94
$ws.Range("I4").Value = 6.4
$ws.Range("L4").Value = 1.49
$ws.Range("V4").Value = 1.19
$ws.Range("Z4").Value = 48
$ws.Range("AB4").Value = 7.8
$ws.Range("AC4").Value = 9.199999999999999
$ws.Range("AH4").Value = 28
$ws.Range("AK4").Value = 25
$ws.Range("AL4").Value = 55
$ws.Range("AN4").Value = 18.5
$ws.Range("F5").Value = 2.58
$ws.Range("G5").Value = 2.9
$ws.Range("H5").Value = 2.86
$ws.Range("I5").Value = 3.3
$ws.Range("J5").Value = 3.05
$ws.Range("K5").Value = 3.5
$ws.Range("O5").Value = 1.38
$ws.Range("P5").Value = 1.71
$ws.Range("Q5").Value = 2.14
$ws.Range("S5").Value = 4.2
$ws.Range("T5").Value = 1.84
$ws.Range("U5").Value = 1.98
$ws.Range("V5").Value = 1.43
$ws.Range("W5").Value = 1.52
$ws.Range("X5").Value = 25
$ws.Range("Y5").Value = 980
$ws.Range("Z5").Value = 1000
$ws.Range("AA5").Value = 1000
$ws.Range("AC5").Value = 7.8
$ws.Range("AD5").Value = 1000
$ws.Range("AF5").Value = 1000
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 60
$ws.Range("AJ5").Value = 980
$ws.Range("G6").Value = 270
$ws.Range("M6").Value = 1.02
$ws.Range("O6").Value = 1.02
$ws.Range("Q6").Value = 1.02
$ws.Range("AC6").Value = 100
$ws.Range("H7").Value = 1.84
$ws.Range("I7").Value = 1.85
$ws.Range("J7").Value = 3.55
$ws.Range("K7").Value = 3.6
$ws.Range("O7").Value = 1.4
$ws.Range("P7").Value = 1.8
$ws.Range("T7").Value = 2.08
$ws.Range("U7").Value = 1.91
$ws.Range("V7").Value = 2.16
$ws.Range("AB7").Value = 16.5
$ws.Range("AK7").Value = 85
$ws.Range("AL7").Value = 90
$ws.Range("AN7").Value = 150
$ws.Range("F8").Value = 2.72
$ws.Range("G8").Value = 3.1
$ws.Range("H8").Value = 2.46
$ws.Range("I8").Value = 2.76
$ws.Range("J8").Value = 3.45
$ws.Range("K8").Value = 3.95
$ws.Range("O8").Value = 1.04
$ws.Range("P8").Value = 1.9
$ws.Range("V8").Value = 1.56
$ws.Range("W8").Value = 1.48
$ws.Range("U9").Value = 1.77
$ws.Range("Y9").Value = 6.4
$ws.Range("F10").Value = 5.6
$ws.Range("G10").Value = 5.7
$ws.Range("H10").Value = 1.75
$ws.Range("I10").Value = 1.76
$ws.Range("J10").Value = 3.9
$ws.Range("K10").Value = 3.95
$ws.Range("N10").Value = 3.8
$ws.Range("O10").Value = 1.34
$ws.Range("P10").Value = 1.96
$ws.Range("Q10").Value = 2.02
$ws.Range("S10").Value = 3.6
$ws.Range("U10").Value = 2.02
$ws.Range("Y10").Value = 8
$ws.Range("Z10").Value = 10.5
$ws.Range("AD10").Value = 9.6
$ws.Range("AF10").Value = 42
$ws.Range("AL10").Value = 80
$ws.Range("J11").Value = 3.5
$ws.Range("K11").Value = 3.55
$ws.Range("N11").Value = 3.6
$ws.Range("V11").Value = 1.59
$ws.Range("X11").Value = 13.5
$ws.Range("AC11").Value = 7.6
$ws.Range("AE11").Value = 75
$ws.Range("AI11").Value = 44
$ws.Range("AO11").Value = 26
$ws.Range("J12").Value = 8.800000000000001
$ws.Range("S12").Value = 1.68
$ws.Range("W12").Value = 6.8
$ws.Range("AD12").Value = 1000
$ws.Range("H13").Value = 2.38
$ws.Range("I13").Value = 2.4
$ws.Range("Q13").Value = 2.02
$ws.Range("R13").Value = 1.37
$ws.Range("V13").Value = 1.71
$ws.Range("G14").Value = 1.58
$ws.Range("H14").Value = 7.4
$ws.Range("K14").Value = 4.4
$ws.Range("N14").Value = 5
$ws.Range("P14").Value = 2.3
$ws.Range("R14").Value = 1.53
$ws.Range("W14").Value = 2.7
$ws.Range("Y14").Value = 28
$ws.Range("Z14").Value = 70
$ws.Range("AE14").Value = 95
$ws.Range("AH14").Value = 19.5
$ws.Range("AM14").Value = 95
$ws.Range("AN14").Value = 6.8

$wb.Save()
